$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("registerSubscriptions")
$ws2 = $wb.Worksheets.Item("delSubscriptionById")

# Update response data for test-id 'iEMS-sub-mgmt-Tes-10' (delSubscriptionById, row 2)
$ws2.Range("I2").Value = "unRegister.clientId: must match"
$ws2.Range("H2").Value = 101400

# Update response data for test-id 'iEMS-sub-mgmt-Test-4' (registerSubscriptions, row 5)
$ws1.Range("B5").Value = "bad request, clientId contains special characters"
$ws1.Range("I5").Value = "register.clientId: must match"
$ws1.Range("G5").Value = 101400

# Update the active sheet/selection to match the edited state
$ws2.Activate()
$ws2.Range("I8").Select()

$ws1.Activate()
$ws1.Range("I5").Select()
